$wb = $excel.ActiveWorkbook

# Sheet "User" (first sheet)
$wsUser = $wb.Worksheets.Item("User")
# D2: "get set" -> "are you set"
$wsUser.Range("D2").Value = "are you set"
# Update selection on User sheet to C27 (B2 text itself, "are you ok", is unchanged)
$wsUser.Range("C27").Select()

# Sheet "Replies" (second sheet)
$wsReplies = $wb.Worksheets.Item("Replies")
# A2: "I am fine,thank you sir" -> "I am fine thank you sir"
$wsReplies.Range("A2").Value = "I am fine thank you sir"
# B2: "yes sir,I am fine" -> "yes sir I am fine"
$wsReplies.Range("B2").Value = "yes sir I am fine"
# D2: "I am ready sir" -> "I am rgood sir"
$wsReplies.Range("D2").Value = "I am rgood sir"

# Replies is the active/selected sheet; set its selection to C2
$wsReplies.Activate()
$wsReplies.Range("C2").Select()
